$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1429913  # H17: was 1429484.2
$ws.Cells.Item(17, 9).Value = 4000  # I17: was 0
$ws.Cells.Item(17, 10).Value = 1667565.1  # J17: was 1429484.2
$ws.Cells.Item(17, 11).Value = 12000  # K17: was 0
$ws.Cells.Item(17, 12).Value = 5002695.300000001  # L17: was 4288452.6
$ws.Cells.Item(17, 13).Value = -11832  # M17: was None
$ws.Cells.Item(17, 14).Value = -5003031.300000001  # N17: was -4288788.6
$ws.Cells.Item(33, 8).Value = 22727892  # H33: was 25000654
$ws.Cells.Item(33, 9).Value = 596.8333  # I33: was 664.4
$ws.Cells.Item(33, 11).Value = 596.8333  # K33: was 664.4
$ws.Cells.Item(33, 13).Value = -367.8333  # M33: was -435.4
$ws.Cells.Item(40, 8).Value = 3432.6667  # H40: was 3527.0908
$ws.Cells.Item(40, 10).Value = 3488.4443  # J40: was 3625.25
$ws.Cells.Item(40, 12).Value = 3488.4443  # L40: was 3625.25
$ws.Cells.Item(40, 14).Value = -3838.4443  # N40: was -3975.25
$ws.Cells.Item(63, 8).Value = 1000000000  # H63: was 0
$ws.Cells.Item(63, 10).Value = 1000000000  # J63: was 0
$ws.Cells.Item(63, 12).Value = 1000000000  # L63: was 0
$ws.Cells.Item(63, 14).Value = -1000001248  # N63: was None
$ws.Cells.Item(64, 8).Value = 4531.25  # H64: was 4277.273
$ws.Cells.Item(64, 9).Value = 4166.6665  # I64: was 3860
$ws.Cells.Item(64, 10).Value = 4750  # J64: was 4625
$ws.Cells.Item(64, 11).Value = 4166.6665  # K64: was 3860
$ws.Cells.Item(64, 12).Value = 4750  # L64: was 4625
$ws.Cells.Item(64, 13).Value = -3918.6665  # M64: was -3612
$ws.Cells.Item(64, 14).Value = -5246  # N64: was -5121
$ws.Cells.Item(66, 8).Value = 1000000000  # H66: was 0
$ws.Cells.Item(66, 10).Value = 1000000000  # J66: was 0
$ws.Cells.Item(66, 12).Value = 3000000000  # L66: was 0
$ws.Cells.Item(66, 14).Value = -3000006240  # N66: was None
$ws.Cells.Item(67, 8).Value = 4531.25  # H67: was 4277.273
$ws.Cells.Item(67, 9).Value = 4166.6665  # I67: was 3860
$ws.Cells.Item(67, 10).Value = 4750  # J67: was 4625
$ws.Cells.Item(67, 11).Value = 4166.6665  # K67: was 3860
$ws.Cells.Item(67, 12).Value = 4750  # L67: was 4625
$ws.Cells.Item(67, 13).Value = -3308.6665  # M67: was -3002
$ws.Cells.Item(67, 14).Value = -6466  # N67: was -6341
$ws.Cells.Item(88, 8).Value = 1362.6923  # H88: was 1395.0834
$ws.Cells.Item(88, 10).Value = 1380.875  # J88: was 1439
$ws.Cells.Item(88, 12).Value = 1380.875  # L88: was 1439
$ws.Cells.Item(88, 14).Value = -2192.875  # N88: was -2251
$ws.Cells.Item(91, 8).Value = 1362.6923  # H91: was 1395.0834
$ws.Cells.Item(91, 10).Value = 1380.875  # J91: was 1439
$ws.Cells.Item(91, 12).Value = 1380.875  # L91: was 1439
$ws.Cells.Item(91, 14).Value = -4188.875  # N91: was -4247
$ws.Cells.Item(100, 8).Value = 3708.0667  # H100: was 4001.3076
$ws.Cells.Item(100, 9).Value = 3943.25  # I100: was 4371.5
$ws.Cells.Item(100, 11).Value = 3943.25  # K100: was 4371.5
$ws.Cells.Item(100, 13).Value = -3402.25  # M100: was -3830.5
$ws.Cells.Item(113, 8).Value = 166670160  # H113: was 66669868
$ws.Cells.Item(113, 10).Value = 4200  # J113: was 3428.5715
$ws.Cells.Item(113, 12).Value = 4200  # L113: was 3428.5715
$ws.Cells.Item(113, 14).Value = -10708  # N113: was -9936.5715
$ws.Cells.Item(116, 8).Value = 36532036  # H116: was 34383150
$ws.Cells.Item(116, 9).Value = 35877244  # I116: was 31392712
$ws.Cells.Item(116, 11).Value = 35877244  # K116: was 31392712
$ws.Cells.Item(116, 13).Value = -35873802  # M116: was -31389270
$ws.Cells.Item(132, 8).Value = 2889.6858  # H132: was 2855.5916
$ws.Cells.Item(132, 9).Value = 2760.07  # I132: was 2720.5688
$ws.Cells.Item(132, 11).Value = 8280.210000000001  # K132: was 8161.7064
$ws.Cells.Item(132, 13).Value = -5750.210000000001  # M132: was -5631.7064

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 110453.7  # H32: was 110454.664
$ws.Cells.Item(32, 9).Value = 118662.95  # I32: was 118663.984
$ws.Cells.Item(32, 11).Value = 118662.95  # K32: was 118663.984
$ws.Cells.Item(32, 13).Value = -118375.95  # M32: was -118376.984
$ws.Cells.Item(122, 8).Value = 3213.7273  # H122: was 3216.7878
$ws.Cells.Item(122, 9).Value = 1632.238  # I122: was 1637.0476
$ws.Cells.Item(122, 11).Value = 4896.714  # K122: was 4911.142800000001
$ws.Cells.Item(122, 13).Value = -2446.714  # M122: was -2461.142800000001
$ws.Cells.Item(132, 8).Value = 910533.75  # H132: was 834741.9399999999
$ws.Cells.Item(132, 9).Value = 527540.9399999999  # I132: was 477397.25
$ws.Cells.Item(132, 11).Value = 1582622.82  # K132: was 1432191.75
$ws.Cells.Item(132, 13).Value = -1580092.82  # M132: was -1429661.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2057.7026  # H20: was 2055.3948
$ws.Cells.Item(20, 9).Value = 2128.5454  # I20: was 2116.4348
$ws.Cells.Item(20, 10).Value = 1953.8  # J20: was 1961.8
$ws.Cells.Item(20, 11).Value = 2128.5454  # K20: was 2116.4348
$ws.Cells.Item(20, 12).Value = 1953.8  # L20: was 1961.8
$ws.Cells.Item(20, 13).Value = -1881.5454  # M20: was -1869.4348
$ws.Cells.Item(20, 14).Value = -2447.8  # N20: was -2455.8
$ws.Cells.Item(80, 8).Value = 83335620  # H80: was 100002540
$ws.Cells.Item(80, 10).Value = 100001780  # J80: was 125001976
$ws.Cells.Item(80, 12).Value = 100001780  # L80: was 125001976
$ws.Cells.Item(80, 14).Value = -100003776  # N80: was -125003972
$ws.Cells.Item(83, 8).Value = 83335620  # H83: was 100002540
$ws.Cells.Item(83, 10).Value = 100001780  # J83: was 125001976
$ws.Cells.Item(83, 12).Value = 500008900  # L83: was 625009880
$ws.Cells.Item(83, 14).Value = -500018884  # N83: was -625019864
$ws.Cells.Item(86, 8).Value = 3000  # H86: was 3152.6667
$ws.Cells.Item(86, 9).Value = 3000  # I86: was 2666.3333
$ws.Cells.Item(86, 10).Value = 0  # J86: was 3639
$ws.Cells.Item(86, 11).Value = 3000  # K86: was 2666.3333
$ws.Cells.Item(86, 12).Value = 0  # L86: was 3639
$ws.Cells.Item(86, 13).Value = -1877  # M86: was -1543.3333
$ws.Cells.Item(86, 14).ClearContents()  # N86: was -5885
$ws.Cells.Item(89, 8).Value = 3000  # H89: was 3152.6667
$ws.Cells.Item(89, 9).Value = 3000  # I89: was 2666.3333
$ws.Cells.Item(89, 10).Value = 0  # J89: was 3639
$ws.Cells.Item(89, 11).Value = 15000  # K89: was 13331.6665
$ws.Cells.Item(89, 12).Value = 0  # L89: was 18195
$ws.Cells.Item(89, 13).Value = -9384  # M89: was -7715.666499999999
$ws.Cells.Item(89, 14).ClearContents()  # N89: was -29427
$ws.Cells.Item(96, 8).Value = 8062.5  # H96: was 7862.5
$ws.Cells.Item(96, 9).Value = 8062.5  # I96: was 7862.5
$ws.Cells.Item(96, 11).Value = 8062.5  # K96: was 7862.5
$ws.Cells.Item(96, 13).Value = -5316.5  # M96: was -5116.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1708.5333  # H22: was 1994
$ws.Cells.Item(22, 10).Value = 2154  # J22: was 2749.25
$ws.Cells.Item(22, 12).Value = 2154  # L22: was 2749.25
$ws.Cells.Item(22, 14).Value = -2854  # N22: was -3449.25
$ws.Cells.Item(31, 8).Value = 2678.94  # H31: was 2660.7273
$ws.Cells.Item(31, 9).Value = 837.9091  # I31: was 838.7273
$ws.Cells.Item(31, 10).Value = 2906.4832  # J31: was 2888.4773
$ws.Cells.Item(31, 11).Value = 837.9091  # K31: was 838.7273
$ws.Cells.Item(31, 12).Value = 2906.4832  # L31: was 2888.4773
$ws.Cells.Item(31, 13).Value = -542.9091  # M31: was -543.7273
$ws.Cells.Item(31, 14).Value = -3496.4832  # N31: was -3478.4773
$ws.Cells.Item(34, 8).Value = 2678.94  # H34: was 2660.7273
$ws.Cells.Item(34, 9).Value = 837.9091  # I34: was 838.7273
$ws.Cells.Item(34, 10).Value = 2906.4832  # J34: was 2888.4773
$ws.Cells.Item(34, 11).Value = 837.9091  # K34: was 838.7273
$ws.Cells.Item(34, 12).Value = 2906.4832  # L34: was 2888.4773
$ws.Cells.Item(34, 13).Value = -635.9091  # M34: was -636.7273
$ws.Cells.Item(34, 14).Value = -3310.4832  # N34: was -3292.4773
$ws.Cells.Item(122, 8).Value = 2960.7144  # H122: was 2692
$ws.Cells.Item(122, 9).Value = 1821.5454  # I122: was 1786.4166
$ws.Cells.Item(122, 10).Value = 7137.6665  # J122: was 4865.4
$ws.Cells.Item(122, 11).Value = 5464.6362  # K122: was 5359.2498
$ws.Cells.Item(122, 12).Value = 21412.9995  # L122: was 14596.2
$ws.Cells.Item(122, 13).Value = -3014.6362  # M122: was -2909.2498
$ws.Cells.Item(122, 14).Value = -26312.9995  # N122: was -19496.2
$ws.Cells.Item(132, 8).Value = 2088.4138  # H132: was 2109.5862
$ws.Cells.Item(132, 9).Value = 1946.8148  # I132: was 1969.5555
$ws.Cells.Item(132, 11).Value = 5840.4444  # K132: was 5908.666499999999
$ws.Cells.Item(132, 13).Value = -3310.4444  # M132: was -3378.666499999999
$ws.Cells.Item(134, 8).Value = 1893.3273  # H134: was 1877.1754
$ws.Cells.Item(134, 9).Value = 1456.2972  # I134: was 1482
$ws.Cells.Item(134, 10).Value = 2791.6667  # J134: was 2608.25
$ws.Cells.Item(134, 11).Value = 4368.8916  # K134: was 4446
$ws.Cells.Item(134, 12).Value = 8375.000100000001  # L134: was 7824.75
$ws.Cells.Item(134, 13).Value = -1833.8916  # M134: was -1911
$ws.Cells.Item(134, 14).Value = -13445.0001  # N134: was -12894.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 83.22221999999999  # H2: was 76.3
$ws.Cells.Item(2, 10).Value = 84.85714  # J2: was 76
$ws.Cells.Item(2, 12).Value = 509.14284  # L2: was 456
$ws.Cells.Item(2, 14).Value = -735.14284  # N2: was -682
$ws.Cells.Item(92, 8).Value = 420.7143  # H92: was 380.375
$ws.Cells.Item(92, 9).Value = 240.83333  # I92: was 220.42857
$ws.Cells.Item(92, 11).Value = 722.49999  # K92: was 661.28571
$ws.Cells.Item(92, 13).Value = 525.50001  # M92: was 586.71429
$ws.Cells.Item(98, 8).Value = 8138.7  # H98: was 8712.223
$ws.Cells.Item(98, 9).Value = 15854.25  # I98: was 13283.2
$ws.Cells.Item(98, 10).Value = 2995  # J98: was 2998.5
$ws.Cells.Item(98, 11).Value = 47562.75  # K98: was 39849.60000000001
$ws.Cells.Item(98, 12).Value = 8985  # L98: was 8995.5
$ws.Cells.Item(98, 13).Value = -46064.75  # M98: was -38351.60000000001
$ws.Cells.Item(98, 14).Value = -11981  # N98: was -11991.5
$ws.Cells.Item(104, 8).Value = 4000  # H104: was 0
$ws.Cells.Item(104, 10).Value = 4000  # J104: was 0
$ws.Cells.Item(104, 12).Value = 12000  # L104: was 0
$ws.Cells.Item(104, 14).Value = -17242  # N104: was None
$ws.Cells.Item(113, 8).Value = 570.875  # H113: was 581.8
$ws.Cells.Item(113, 9).Value = 355  # I113: was 341
$ws.Cells.Item(113, 10).Value = 620.6923  # J113: was 669.36365
$ws.Cells.Item(113, 11).Value = 1065  # K113: was 1023
$ws.Cells.Item(113, 12).Value = 1862.0769  # L113: was 2008.09095
$ws.Cells.Item(113, 13).Value = 1105  # M113: was 1147
$ws.Cells.Item(113, 14).Value = -6202.0769  # N113: was -6348.09095
$ws.Cells.Item(131, 8).Value = 11145.261  # H131: was 10324.24
$ws.Cells.Item(131, 9).Value = 1666  # I131: was 1504.5
$ws.Cells.Item(131, 10).Value = 13778.389  # J131: was 13109.421
$ws.Cells.Item(131, 11).Value = 4998  # K131: was 4513.5
$ws.Cells.Item(131, 12).Value = 41335.167  # L131: was 39328.263
$ws.Cells.Item(131, 13).Value = 42  # M131: was 526.5
$ws.Cells.Item(131, 14).Value = -51415.167  # N131: was -49408.263
$ws.Cells.Item(138, 8).Value = 5465502.5  # H138: was 6011756
$ws.Cells.Item(138, 10).Value = 22594  # J138: was 27500
$ws.Cells.Item(138, 12).Value = 67782  # L138: was 82500
$ws.Cells.Item(138, 14).Value = -78062  # N138: was -92780

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3476.1333  # H122: was 3236.0667
$ws.Cells.Item(122, 9).Value = 3114.3  # I122: was 2961.75
$ws.Cells.Item(122, 10).Value = 4199.8  # J122: was 4333.3335
$ws.Cells.Item(122, 11).Value = 9342.900000000001  # K122: was 8885.25
$ws.Cells.Item(122, 12).Value = 12599.4  # L122: was 13000.0005
$ws.Cells.Item(122, 13).Value = -6892.900000000001  # M122: was -6435.25
$ws.Cells.Item(122, 14).Value = -17499.4  # N122: was -17900.0005
$ws.Cells.Item(126, 8).Value = 11144.6  # H126: was 9128.429
$ws.Cells.Item(126, 10).Value = 5000  # J126: was 4452.8
$ws.Cells.Item(126, 12).Value = 15000  # L126: was 13358.4
$ws.Cells.Item(126, 14).Value = -19940  # N126: was -18298.4
$ws.Cells.Item(130, 8).Value = 0  # H130: was 57500
$ws.Cells.Item(130, 10).Value = 0  # J130: was 57500
$ws.Cells.Item(130, 12).Value = 0  # L130: was 57500
$ws.Cells.Item(130, 14).ClearContents()  # N130: was -67540
$ws.Cells.Item(132, 8).Value = 671846.9399999999  # H132: was 458358.12
$ws.Cells.Item(132, 9).Value = 1116021.5  # I132: was 628085.5
$ws.Cells.Item(132, 10).Value = 5585.1665  # J132: was 5751.8335
$ws.Cells.Item(132, 11).Value = 3348064.5  # K132: was 1884256.5
$ws.Cells.Item(132, 12).Value = 16755.4995  # L132: was 17255.5005
$ws.Cells.Item(132, 13).Value = -3345534.5  # M132: was -1881726.5
$ws.Cells.Item(132, 14).Value = -21815.4995  # N132: was -22315.5005
$ws.Cells.Item(133, 8).Value = 80000  # H133: was 90435.8
$ws.Cells.Item(133, 9).Value = 80000  # I133: was 84999
$ws.Cells.Item(133, 10).Value = 0  # J133: was 91795
$ws.Cells.Item(133, 11).Value = 80000  # K133: was 84999
$ws.Cells.Item(133, 12).Value = 0  # L133: was 91795
$ws.Cells.Item(133, 13).Value = -74940  # M133: was -79939
$ws.Cells.Item(133, 14).ClearContents()  # N133: was -101915
$ws.Cells.Item(135, 8).Value = 90000  # H135: was 89950
$ws.Cells.Item(135, 10).Value = 90000  # J135: was 89950
$ws.Cells.Item(135, 12).Value = 90000  # L135: was 89950
$ws.Cells.Item(135, 14).Value = -100140  # N135: was -100090

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 7087.0527  # H132: was 6618.364
$ws.Cells.Item(132, 9).Value = 3724.9  # I132: was 3707.6155
$ws.Cells.Item(132, 11).Value = 11174.7  # K132: was 11122.8465
$ws.Cells.Item(132, 13).Value = -8644.700000000001  # M132: was -8592.8465
$ws.Cells.Item(136, 8).Value = 9060.277  # H136: was 9692.875
$ws.Cells.Item(136, 9).Value = 5560.1  # I136: was 5950.25
$ws.Cells.Item(136, 11).Value = 16680.3  # K136: was 17850.75
$ws.Cells.Item(136, 13).Value = -14130.3  # M136: was -15300.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 55627  # H45: was 61090
$ws.Cells.Item(45, 9).Value = 42840.5  # I45: was 69420
$ws.Cells.Item(45, 10).Value = 59280.285  # J45: was 59900
$ws.Cells.Item(45, 11).Value = 42840.5  # K45: was 69420
$ws.Cells.Item(45, 12).Value = 59280.285  # L45: was 59900
$ws.Cells.Item(45, 13).Value = -42349.5  # M45: was -68929
$ws.Cells.Item(45, 14).Value = -60262.285  # N45: was -60882
$ws.Cells.Item(96, 8).Value = 4595.5386  # H96: was 4562.6665
$ws.Cells.Item(96, 10).Value = 5177.8  # J96: was 5224.75
$ws.Cells.Item(96, 12).Value = 5177.8  # L96: was 5224.75
$ws.Cells.Item(96, 14).Value = -7923.8  # N96: was -7970.75
$ws.Cells.Item(107, 8).Value = 376.55  # H107: was 377.15
$ws.Cells.Item(107, 9).Value = 369.89474  # I107: was 370.5263
$ws.Cells.Item(107, 11).Value = 1109.68422  # K107: was 1111.5789
$ws.Cells.Item(107, 13).Value = 810.3157799999999  # M107: was 808.4211
$ws.Cells.Item(122, 8).Value = 5436.8  # H122: was 5822.1113
$ws.Cells.Item(122, 9).Value = 4281.143  # I122: was 4666.5
$ws.Cells.Item(122, 11).Value = 12843.429  # K122: was 13999.5
$ws.Cells.Item(122, 13).Value = -10393.429  # M122: was -11549.5
$ws.Cells.Item(132, 8).Value = 559704.25  # H132: was 479884.25
$ws.Cells.Item(132, 9).Value = 873594.9399999999  # I132: was 693050.7
$ws.Cells.Item(132, 11).Value = 2620784.82  # K132: was 2079152.1
$ws.Cells.Item(132, 13).Value = -2618254.82  # M132: was -2076622.1
$ws.Cells.Item(135, 8).Value = 99999  # H135: was 74999.5
$ws.Cells.Item(135, 10).Value = 0  # J135: was 50000
$ws.Cells.Item(135, 12).Value = 0  # L135: was 50000
$ws.Cells.Item(135, 14).ClearContents()  # N135: was -60140
$ws.Cells.Item(136, 8).Value = 7129  # H136: was 7467.8945
$ws.Cells.Item(136, 9).Value = 5786.625  # I136: was 6126.4
$ws.Cells.Item(136, 11).Value = 17359.875  # K136: was 18379.2
$ws.Cells.Item(136, 13).Value = -14809.875  # M136: was -15829.2
